# Applies the "Updated cryptos list" data refresh described in the commit.
# Only the Price (D) and Volume(1h) (E) columns change; every other cell is
# left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the Price column as Text before writing into it. Several of the new
# price strings (e.g. "609.48") would otherwise be auto-converted by Excel
# into numeric values instead of being kept as the plain text the workbook
# originally stored them as.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "65.866.84"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "2.695.48"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "609.48"
$ws.Range("E5").Value = "  +1.60%  "
$ws.Range("D6").Value = "157.75"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("E9").Value = "  +4.92%  "
$ws.Range("D10").Value = "6.05"
$ws.Range("E10").Value = "  +4.71%  "
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "30.38"
$ws.Range("E13").Value = "  +5.23%  "
$ws.Range("E14").Value = "  +8.82%  "
$ws.Range("D15").Value = "3.178.89"
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").Value = "65.723.10"
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("D17").Value = "2.686.41"
$ws.Range("E17").Value = "  +2.79%  "
$ws.Range("D18").Value = "12.68"
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("D19").Value = "4.90"
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("D20").Value = "360.51"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("D21").Value = "7.60"
$ws.Range("E21").Value = "  +5.17%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "70.21"
$ws.Range("E23").Value = "  +3.35%  "
$ws.Range("D24").Value = "9.84"
$ws.Range("E24").Value = "  +3.97%  "
$ws.Range("D25").Value = "0.0000107"
$ws.Range("E25").Value = "  +12.37%  "
$ws.Range("E26").Value = "  -4.14%  "
$ws.Range("E27").Value = "  +2.78%  "
$ws.Range("D28").Value = "0.169"
$ws.Range("E28").Value = "  +3.62%  "
$ws.Range("D29").Value = "8.27"
$ws.Range("E29").Value = "  +2.45%  "
$ws.Range("E30").Value = "  +4.49%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "531.91"
$ws.Range("E32").Value = "  +2.78%  "
$ws.Range("D33").Value = "1.80"
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("D34").Value = "6.66"
$ws.Range("E34").Value = "  +5.21%  "
$ws.Range("D35").Value = "5.46"
$ws.Range("E35").Value = "  -3.66%  "
$ws.Range("D36").Value = "0.433"
$ws.Range("E36").Value = "  +1.39%  "
$ws.Range("D37").Value = "20.82"
$ws.Range("E37").Value = "  +3.21%  "
$ws.Range("D38").Value = "163.31"
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("E39").Value = "  -2.49%  "
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D42").Value = "169.25"
$ws.Range("E42").Value = "  +2.70%  "
$ws.Range("D43").Value = "42.60"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").Value = "4.18"
$ws.Range("E44").Value = "  +2.28%  "
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("E46").Value = "  +1.64%  "
$ws.Range("D47").Value = "2.27"
$ws.Range("E47").Value = "  +2.80%  "
$ws.Range("D48").Value = "0.659"
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("D49").Value = "0.0266"
$ws.Range("E49").Value = "  +4.49%  "
$ws.Range("D50").Value = "21.27"
$ws.Range("E50").Value = "  +9.24%  "
$ws.Range("D51").Value = "0.0984"
$ws.Range("E51").Value = "  +0.11%  "

# Put the Price column back to the workbook's default (unstyled) look now
# that the text values are safely stored, so no stray explicit style index
# is left on cells that originally had none.
$ws.Range("D2:D51").Style = "Normal"

Write-Host "Updated cryptos list on $(Get-Date) with GitHub Actions"
